$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new data rows to append after the existing last row (114),
# continuing the same index/style pattern used throughout the sheet.
$newRows = @(
    @{ Row = 115; Idx = 113; Code = "703";  Comune = 61 },
    @{ Row = 116; Idx = 114; Code = ".74";  Comune = 134 },
    @{ Row = 117; Idx = 115; Code = ".207"; Comune = 388 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    # Column A: running index, formatted like the rest of column A
    # (bold / bordered / centered-top) by copying the format from the
    # last existing row in that column.
    $ws.Range("A114").Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $item.Idx

    # Column B: particle code stored as plain text (values like "703"
    # or ".74" must stay text, not become numbers). Build the text via a
    # formula on a scratch cell (so Excel's type is definitely "string"),
    # then copy just the resulting value into place - this keeps the
    # cell's number format/style completely untouched (no style churn).
    $ws.Range("Z1").Formula = "=""" + $item.Code + """"
    $ws.Range("Z1").Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4163)

    # Column C: plain numeric value.
    $ws.Cells.Item($r, 3).Value = $item.Comune
}

# Clean up the scratch cell used to mint text values.
$ws.Range("Z1").Clear()
